$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Trim the first paragraph's text:
#    "Daily morning Check Chats, Mails  and then Download latest build for
#     both Platforms " -> "Daily morning Check Chats, Mails  and then
#     Download latest build"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Daily morning Check Chats, Mails  and then Download latest build for both Platforms ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Daily morning Check Chats, Mails  and then Download latest build", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the two runs "2)A" + "ttach Video to Bug/PBI" into a single run
#    of text "2)Attach Video to Bug/PBI"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "2)A" + "ttach Video to Bug/PBI",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2)Attach Video to Bug/PBI", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert two new paragraphs right before the "Before requesting for ..."
#    paragraph (immediately after the "Automation-" paragraph):
#      "Use Perfecto Report Library Report instead of Terminal"
#      <empty paragraph>
# ---------------------------------------------------------------------------
$f = $d.Content
$f.Find.Execute("Automation-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endOfAutomation = $f.End

# Insert paragraph mark right after "Automation-" -> creates a new, empty
# paragraph between the "Automation-" paragraph and "Before requesting...".
$insPt = $d.Range($endOfAutomation, $endOfAutomation)
$insPt.InsertParagraphAfter()

# Fill the new paragraph's text.
$newText = "Use Perfecto Report Library Report instead of Terminal"
$newParaRange = $d.Range($endOfAutomation + 1, $endOfAutomation + 1)
$newParaRange.Text = $newText

# Insert a second, empty paragraph after the one we just filled.
$endOfNewPara = $endOfAutomation + 1 + $newText.Length
$insPt2 = $d.Range($endOfNewPara, $endOfNewPara)
$insPt2.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4) Replace "master Security Token and No commented Spec files" with
#    "delete ios.config.ts" (keeping the preceding "Check " run intact)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "master Security Token and No commented Spec files",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "delete ios.config.ts", 2) | Out-Null
